$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.794.88'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.489.65'
$ws.Range('E3').Value = '  +0.70%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '588.02'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '175.24'
$ws.Range('E6').Value = '  +2.12%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.146'
$ws.Range('E9').Value = '  +5.45%  '
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.337'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.947.04'
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.53'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '67.707.09'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.478.43'
$ws.Range('E17').Value = '  +2.16%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.89'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.46'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '349.81'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.14'
$ws.Range('E21').Value = '  +3.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '71.21'
$ws.Range('E22').Value = '  +4.26%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.25'
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.71'
$ws.Range('E25').Value = '  -4.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.97'
$ws.Range('E26').Value = '  -2.94%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.582.41'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0899'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.77'
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '497.71'
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.76'
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '164.51'
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('E36').Value = '  +3.21%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.63'
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.33'
$ws.Range('E38').Value = '  +1.04%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.31'
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.73'
$ws.Range('E41').Value = '  +2.74%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.328'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.81'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.44'
$ws.Range('E44').Value = '  +3.62%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '147.22'
$ws.Range('E45').Value = '  +3.13%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.55'
$ws.Range('E46').Value = '  +2.45%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0₆0258'
$ws.Range('E47').Value = '  +2.73%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.514'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('E49').Value = '  +0.91%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.57'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.580'
$ws.Range('E51').Value = '  -0.55%  '
